$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (F1) onto the two new
# header cells so they get the same bold/border/centered formatting (style
# index 1), then set their text.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# New data cells for row 2 (no special style, same as the other data cells)
$ws.Range("G2").Value = 0.1180509527000443
$ws.Range("H2").Value = 0.991
